$d = $word.ActiveDocument

# --- 1. Simplify Index entries (merge multi-run paragraphs into single runs) ---
$null = $d.Content.Find.Execute("3. Tech Stack", $false, $false, $false, $false, $false, $true, 1, $false, "3. Tech Stack", 2)
$null = $d.Content.Find.Execute("4. Folder Structure", $false, $false, $false, $false, $false, $true, 1, $false, "4. Folder Structure", 2)
$null = $d.Content.Find.Execute("5. Logic", $false, $false, $false, $false, $false, $true, 1, $false, "5. Logic", 2)

# --- 2. Merge "UserScore model." run split (removes spell-check markup, joins runs) ---
$oldModelText = " model. The creativity assessments which have only a few questions or less data are not using the database to fetch questions. Instead, I store the questions in the resources/json "
$null = $d.Content.Find.Execute($oldModelText, $false, $false, $false, $false, $false, $true, 1, $false, $oldModelText, 2)

# --- 3. Insert new "app pages" details block + Assessment Logic section at end of document ---
$lastPara = $d.Paragraphs.Last.Range
$newXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Home page: Lists the assessments by category. Each assessment contains a short description and an Open button. The button takes the user to the Assessment Details page (/assessments/:key)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Assessment Details Page: This page contains writeup on assessment. It has a start button to start the assessment. If the user is already enrolled in the assessment they are taken to the Assessment page. Otherwise, they are first enrolled and then redirected to the Assessment page.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Assessment Page: Except the psychometric assessments, the assessment page is a single page containing the questionnaire and inputs. The user fills out the form/ enters their answers. There will be a button to Save and Submit. The save button will save the user responses in the database. The user responses will be populated on page load. The user can continue where they left off. The submit button will submit the assessment. This will also save their responses as well as grade the assessment with the responses. Each assessment controls the scoring individually. For NEST and CPT there is an additional page. The first page lists the different modules or sections of the assessment. The user can choose to answer any of the sections first and save their responses at each step. The assessment page will also show statistics about the modules such as the no of questions they have completed and time spent so far. On submitting the entire form the score is created at the backend and the user is taken to the scores page.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Score page: The score page fetches all the past user scores and shows them on the page. It also shows the details about the assessment such as the time when the assessment was submitted, a graph showing the score for each section of the assessment and a table showing the data in the graph. For each user score, they have the option to retake the test, delete the score or view it’s report. Retaking the test will take them back to the assessment page, where there previous responses are still saved. The user scores are retained for every attempt, and a new score is simply created each time the user hits submit.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Profile page: The profile page contains the user information and the assessment list. The assessment list contains the assessments the user is enrolled in. They have an option to open the assessment or remove it. On opening the assessment, they are taken back to the assessment page. On removing, the user assessment is deleted and does not show on their profile anymore. The user scores are not deleted. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Auth: Login page</w:t></w:r><w:r><w:t>: The user can login with their email and password, or with google.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Sign</w:t></w:r><w:r><w:t xml:space="preserve"> up page</w:t></w:r><w:r><w:t>: The user registers their email and password to create a new profile. Then they can login with their credentials.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>P</w:t></w:r><w:r><w:t>rofile update page</w:t></w:r><w:r><w:t xml:space="preserve">: Here the user can update personal information about their profile. This information is also reflected in the reports. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Report page: If there is a report available for the assessment they will be shown to the user. The report will contain the user information, their scores, graph, and additionally a write up on how to interpret their score. </w:t></w:r><w:r><w:t xml:space="preserve">The user can download the graph as an image from either the score page or the report page. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Divergent Responses page: This is for the administrator to accept and reject the responses for the divergent thinking assessment. It is the only test for which the user does not immediately see their result. The result is shown after the administrator publishes their result. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">All scores will appear automatically on the score page. The user can login into the website to check them at any time. The reports are dynamically created and can be viewed at any time as well. They can print the page to obtain the report. The assessments can be taken multiple times. The scores can also be deleted and recreated as they wish. The reports for some assessment will contain static data, however for some assessment the report varies from person to person. The user is shown feedback depending on their score. </w:t></w:r></w:p><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t>Assessment Logic</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>NEST assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>CPT assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Creativity assessments</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Creativity personality assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Creativity environment assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Creativity motivation assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Creativity temperament assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Creativity convergent thinking assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Creativity divergent thinking assessment</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Social leadership assessment</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $lastPara.InsertXML($newXml)
